$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")
$ws.Activate()

# Add new "Alarm Normal load method" detail columns (S, T) next to the
# existing Slot Card detail block (I..M / N..R), mirroring the header
# and data row styling used by the other header/data cells.

# Header row (row 7): labels for the new detail columns
$ws.Range("S7").Value = "AlarmLoadingDetail"
$ws.Range("T7").Value = "StandbyLoadingDetail"
$ws.Range("A7").Copy()
$ws.Range("S7:T7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row (row 8): values for the new detail columns
$ws.Range("S8").Value = "Battery Alarm (A)"
$ws.Range("T8").Value = "Battery Standby (A)"
$ws.Range("A8").Copy()
$ws.Range("S8:T8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update page setup (matches the other worksheet in the workbook)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Update the visible selection / scroll position to the newly added
# columns, as left by the editor that made this change.
$excel.ActiveWindow.ScrollColumn = 4
$null = $ws.Range("S7:T8").Select()
